# Updates cryptos list with latest price/volume data (GitHub Actions scheduled refresh).
# Source diff only touches columns B-E of rows 2-51 on the single worksheet;
# column A (rank index) and row 1 (header) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value, in the same row order as the source diff.
$updates = [ordered]@{
    'D2' = '96.192.65'
    'E2' = '  -0.01%  '
    'D3' = '3.521.37'
    'E3' = '  +6.22%  '
    'D4' = '1.00'
    'E4' = '  -0.04%  '
    'D5' = '243.60'
    'E5' = '  -1.19%  '
    'D6' = '652.48'
    'E6' = '  +0.39%  '
    'D7' = '1.45'
    'E7' = '  +7.17%  '
    'D8' = '0.419'
    'E8' = '  +1.31%  '
    'D9' = '1.00'
    'E9' = '  +0.03%  '
    'D10' = '1.00'
    'E10' = '  +2.46%  '
    'D11' = '3.521.51'
    'E11' = '  +6.28%  '
    'D12' = '43.29'
    'E12' = '  +8.90%  '
    'E13' = '  -1.55%  '
    'D14' = '6.20'
    'E14' = '  +3.70%  '
    'D15' = '95.932.57'
    'E15' = '  +0.02%  '
    'D16' = '4.190.60'
    'E16' = '  +6.67%  '
    'D17' = '0.0000256'
    'E17' = '  +2.98%  '
    'D18' = '8.59'
    'E18' = '  +1.12%  '
    'D19' = '3.540.48'
    'E19' = '  +6.69%  '
    'D20' = '18.59'
    'E20' = '  +11.64%  '
    'D21' = '12.21'
    'E21' = '  +18.16%  '
    'D22' = '0.510'
    'E22' = '  +7.33%  '
    'D23' = '517.44'
    'E23' = '  +4.44%  '
    'E24' = '  -0.03%  '
    'D25' = '0.0000197'
    'E25' = '  +0.79%  '
    'E26' = '  +4.65%  '
    'D27' = '92.84'
    'E27' = '  -1.09%  '
    'D28' = '12.60'
    'E28' = '  +6.22%  '
    'D29' = '3.712.52'
    'E29' = '  +6.48%  '
    'D30' = '12.06'
    'E30' = '  +12.81%  '
    'E31' = '  +15.36%  '
    'D32' = '1.00'
    'E32' = '  +0.12%  '
    'D33' = '0.140'
    'E33' = '  -0.85%  '
    'E34' = '  +1.40%  '
    'D35' = '31.72'
    'E35' = '  +14.16%  '
    'D36' = '0.586'
    'E36' = '  +8.89%  '
    'D37' = '0.998'
    'E37' = '  -0.40%  '
    'D38' = '7.99'
    'E38' = '  +7.32%  '
    'D39' = '1.50'
    'E39' = '  +2.48%  '
    'B40' = 'Bittensor'
    'C40' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D40' = '525.30'
    'E40' = '  +4.87%  '
    'B41' = 'Kaspa'
    'C41' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D41' = '0.153'
    'E41' = '  +2.66%  '
    'B42' = 'ARBITRUM'
    'C42' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D42' = '0.937'
    'E42' = '  +14.71%  '
    'E43' = '  +0.01%  '
    'D44' = '24.17'
    'E44' = '  -1.12%  '
    'E45' = '  +9.32%  '
    'D46' = '0.0426'
    'E46' = '  +6.59%  '
    'D47' = '3.64'
    'E47' = '  -0.90%  '
    'D48' = '5.63'
    'E48' = '  +3.56%  '
    'D49' = '3.34'
    'E49' = '  +6.99%  '
    'D50' = '2.21'
    'E50' = '  +13.76%  '
    'D51' = '8.33'
    'E51' = '  +0.45%  '
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)

    if ($cellRef -match '^[DE]') {
        # Price/volume columns hold plain text (e.g. "96.192.65", "1.00", "  +0.42%  ").
        # Force text format first so Excel doesn't reinterpret numeric-looking
        # strings (trailing zeros, thousand-dot separators) as real numbers,
        # then restore the default "Normal" style so formatting stays untouched.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    } else {
        # Coin name / link columns are ordinary text already.
        $cell.Value = $newValue
    }
}
